$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Move the selection on the first sheet (Tabelle1) and drop its "active" tab state.
$ws1.Range("D50").Select() | Out-Null

# Add a new worksheet right after Tabelle1 and make it the active sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Col box mod"

# --- Header row ---
$ws2.Range("C6").Value = "Golovin:"
$ws2.Range("E6").Value = "Hall"
$ws2.Range("H6").Value = "Long"

# --- kappa row ---
$ws2.Range("C7").Value = "kappa"
$ws2.Range("E7").Value = "kappa"
$ws2.Range("H7").Value = "kappa"

# --- numeric grid ---
$ws2.Range("C8").Value = 5
$ws2.Range("E8").Value = 5
$ws2.Range("H8").Value = 5

$ws2.Range("C9").Value = 10
$ws2.Range("E9").Value = 10
$ws2.Range("H9").Value = 10

$ws2.Range("C10").Value = 20
$ws2.Range("E10").Value = 20
$ws2.Range("H10").Value = 20

$ws2.Range("C11").Value = 40
$ws2.Range("E11").Value = 40
$ws2.Range("H11").Value = 40

$ws2.Range("C12").Value = 60
$ws2.Range("E12").Value = 100
$ws2.Range("H12").Value = 100

$ws2.Range("C13").Value = 100
$ws2.Range("E13").Value = 200
$ws2.Range("H13").Value = 200

$ws2.Range("E14").Value = 400
$ws2.Range("F14").Value = "finished"
$ws2.Range("H14").Value = 400
$ws2.Range("I14").Value = "finished"

$ws2.Range("E15").Value = 1000
$ws2.Range("F15").Value = "only to seed 1243"
$ws2.Range("H15").Value = 1000
$ws2.Range("I15").Value = "only to seed 1227"

$ws2.Range("E16").Value = 2000
$ws2.Range("H16").Value = 2000

$ws2.Range("E17").Value = 3000
$ws2.Range("H17").Value = 3000

# --- column widths --- (columns A-E and G onward already sit at the sheet's
# default width, which matches the 10.5 used throughout; only F needs a
# wider, explicit override)
$ws2.Columns.Item(6).ColumnWidth = 19.42

# --- selection / active cell on the new sheet ---
$ws2.Range("G30").Select() | Out-Null
